$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 2

$ws.Range("F1").Select()
